$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 24).Value = 1
    $ws.Cells.Item($row, 25).Value = "x"
}
